$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'59.138.27"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.05%  "
$ws.Range("D3").Value = "'2.588.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.36%  "
$ws.Range("E4").Value = "  +0.20%  "
$ws.Range("D5").Value = "'525.00"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").Value = "'139.63"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.95%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  -0.94%  "
$ws.Range("D9").Value = "'2.601.75"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("D10").Value = "'6.51"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("E11").Value = "  -0.29%  "
$ws.Range("D12").Value = "'0.331"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.56%  "
$ws.Range("E13").Value = "  +2.91%  "
$ws.Range("D14").Value = "'3.048.97"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "'59.103.13"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("D16").Value = "'20.55"
$ws.Range("D16").Style = "Normal"
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "'2.626.44"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.09%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "'0.0000133"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.67%  "
$ws.Range("D19").Value = "'340.93"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.08%  "
$ws.Range("E20").Value = "  -0.99%  "
$ws.Range("D21").Value = "'10.11"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.58%  "
$ws.Range("D22").Value = "'6.43"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.10%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'66.61"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.66%  "
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("D26").Value = "'0.405"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.32%  "
$ws.Range("D27").Value = "'0.997"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.18%  "
$ws.Range("D28").Value = "'7.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.51%  "
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").Value = "'0.0₃0726"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.23%  "
$ws.Range("D31").Value = "'5.96"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.96%  "
$ws.Range("E32").Value = "  +0.46%  "
$ws.Range("D33").Value = "'18.72"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.53%  "
$ws.Range("D34").Value = "'149.12"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.43%  "
$ws.Range("D35").Value = "'3.98"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -1.00%  "
$ws.Range("E36").Value = "  -1.23%  "
$ws.Range("D37").Value = "'36.81"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.16%  "
$ws.Range("D38").Value = "'1.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.65%  "
$ws.Range("D39").Value = "'0.830"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.00%  "
$ws.Range("D40").Value = "'0.810"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -7.18%  "
$ws.Range("E41").Value = "  -0.37%  "
$ws.Range("D42").Value = "'0.998"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.00%  "
$ws.Range("D43").Value = "'272.49"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.34%  "
$ws.Range("D44").Value = "'0.603"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("E45").Value = "  +0.85%  "
$ws.Range("E46").Value = "  -0.72%  "
$ws.Range("D47").Value = "'0.0515"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.74%  "
$ws.Range("D48").Value = "'18.47"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.94%  "
$ws.Range("D49").Value = "'1.968.46"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.29%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0223"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.30%  "
$ws.Range("B51").Value = "InjectiveProtocol"
$ws.Range("C51").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D51").Value = "'18.31"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.10%  "
